# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date, Publisher values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/enrollment-pcp-name-on-enrollment"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": clear the stray Constraint(s) value on the root
#     "Extension" row (row 2, column AI) -- it belongs only on the
#     Extension.extension child row (row 4), which already has it. ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The "Fixed Value" for Extension.url (row 5, column Q) mirrors the
# StructureDefinition URL on the Metadata sheet, so it must pick up the
# same rename.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/enrollment-pcp-name-on-enrollment"
